# Latest results and figures
# Update the "data-superpg" sheet's first benchmark table (rows 7-16):
# rows 7-15 get new counter labels ("...,usr" / "(os + usr)" variants) and
# refreshed measured values; row 16 (previously "elapse time") is cleared out
# entirely since that series moved up into row 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data-superpg")

# --- Row 7: CPU_CLK_UNHALTED.THREAD_P -> CPU_CLK_UNHALTED.THREAD_P,usr ---
$ws.Range("A7").Value = "CPU_CLK_UNHALTED.THREAD_P,usr"
$ws.Range("B7").Value = 269553674.51539999
$ws.Range("C7").Value = 270613181.92799997
$ws.Range("D7").Value = 254968323.70770001
$ws.Range("E7").Value = 254986394.25944999

# --- Row 8: ITLB_MISSES.WALK_COMPLETED -> DTLB_LOAD_MISSES.WALK_PENDING,usr ---
$ws.Range("A8").Value = "DTLB_LOAD_MISSES.WALK_PENDING,usr"
$ws.Range("B8").Value = 1587749.1723499999
$ws.Range("C8").Value = 1623758.652
$ws.Range("D8").Value = 657411.55735000002
$ws.Range("E8").Value = 647626.05870000005

# --- Row 9: ITLB_MISSES.WALK_PENDING -> DTLB_STORE_MISSES.WALK_PENDING,usr ---
$ws.Range("A9").Value = "DTLB_STORE_MISSES.WALK_PENDING,usr"
$ws.Range("B9").Value = 539639.92255000002
$ws.Range("C9").Value = 549522.79244999995
$ws.Range("D9").Value = 435870.00605000003
$ws.Range("E9").Value = 418838.23664999998

# --- Row 10: ICACHE_64B.IFTAG_STALL -> ITLB_MISSES.WALK_PENDING,usr ---
$ws.Range("A10").Value = "ITLB_MISSES.WALK_PENDING,usr"
$ws.Range("B10").Value = 2683057.7054499998
$ws.Range("C10").Value = 2691092.7038500002
$ws.Range("D10").Value = 193876.04384999999
$ws.Range("E10").Value = 148342.50774999999

# --- Row 11: CPU_CLK_UNHALTED.THREAD_P (os + usr) -> ICACHE_64B.IFTAG_STALL,usr ---
$ws.Range("A11").Value = "ICACHE_64B.IFTAG_STALL,usr"
$ws.Range("B11").Value = 14934563.44795
$ws.Range("C11").Value = 15033884.34285
$ws.Range("D11").Value = 5605968.9563999996
$ws.Range("E11").Value = 5555615.6648000004

# --- Row 12: DTLB_LOAD_MISSES.WALK_COMPLETED -> CPU_CLK_UNHALTED.THREAD_P (os + usr) ---
$ws.Range("A12").Value = "CPU_CLK_UNHALTED.THREAD_P (os + usr)"
$ws.Range("B12").Value = 298525048.25484997
$ws.Range("C12").Value = 299636984.59315002
$ws.Range("D12").Value = 278691753.02929997
$ws.Range("E12").Value = 278672485.94489998

# --- Row 13: DTLB_LOAD_MISSES.WALK_PENDING -> INST_RETIRED.ANY_P (os + usr) ---
$ws.Range("A13").Value = "INST_RETIRED.ANY_P (os + usr)"
$ws.Range("B13").Value = 249681812.76295
$ws.Range("C13").Value = 249649361.74445
$ws.Range("D13").Value = 239455846.833
$ws.Range("E13").Value = 239755297.35315001

# --- Row 14: DTLB_STORE_MISSES.WALK_COMPLETED -> INST_RETIRED.ANY_P,usr ---
$ws.Range("A14").Value = "INST_RETIRED.ANY_P,usr"
$ws.Range("B14").Value = 223635326.6234
$ws.Range("C14").Value = 223635313.88734999
$ws.Range("D14").Value = 223635302.04449999
$ws.Range("E14").Value = 223635180.03619999

# --- Row 15: DTLB_STORE_MISSES.WALK_PENDING -> elaspe time ---
$ws.Range("A15").Value = "elaspe time"
$ws.Range("B15").Value = 434.46350000000001
$ws.Range("C15").Value = 436.22449999999998
$ws.Range("D15").Value = 400.39800000000002
$ws.Range("E15").Value = 400.3895

# Rows 7-15 lose their red-font ("DTLB/ITLB walk" highlight) styling -
# they revert to the sheet's default (unstyled) look.
$ws.Range("A7:A15").Style = "Normal"

# --- Row 16: previously "elapse time" data row, now emptied entirely ---
$ws.Range("A16:I16").ClearContents()

# Update the current selection to match the saved view.
$ws.Activate()
$ws.Range("B10:E10").Select()
